$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-15 down to 11-16
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's new record
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44413
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100112035
$ws.Cells.Item(10, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 25
$ws.Cells.Item(10, 11).Value = 24000
$ws.Cells.Item(10, 12).Value = 25000
$ws.Cells.Item(10, 13).Value = 24480
$ws.Cells.Item(10, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(10, 15).Value = "Hijuelas"
$ws.Cells.Item(10, 16).Value = 1632
$ws.Cells.Item(10, 17).Value = 15
$ws.Cells.Item(10, 18).Value = "Hortaliza"
